$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.051944486854919064
$ws.Range("B1").Value = -0.051944486884242344

$ws.Range("A2").Value = 0.013464490202881409
$ws.Range("B2").Value = -0.013464490252297064

$ws.Range("A3").Value = -0.020799110086841498
$ws.Range("B3").Value = 0.02079910999601697

$ws.Range("A4").Value = 0.028503384279370134
$ws.Range("B4").Value = -0.028503384303528591
